# NaukriSearch.xlsx update
# 1. Move the current top search result (row 2) on "NaukriSearch" down into the
#    history sheet "Sheet1" as a new last row.
# 2. Write a brand-new top search result ("ServiceNow PM") into row 2 of
#    "NaukriSearch", re-using the same look/formatting as the previous entry.
# 3. Append an additional history row ("Sitecore") to "Sheet1".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NaukriSearch")
$ws2 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Step 1: move NaukriSearch!A2:K2 -> Sheet1!A10:K10 (same values + formats)
# ---------------------------------------------------------------------------
$ws1.Range("A2:K2").Copy() | Out-Null
$ws2.Range("A10:K10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws2.Range("A10").Value2 = 1
$ws2.Range("B10").Value2 = "DotNet Developer"
$ws2.Range("C10").Value2 = "Asp.net, C#, MVC,SQL,Angular"
$ws2.Range("D10").Value2 = "Asp.net Developer,.Net Developer,Full stack.Net Developer,Dotnet Developer"
$ws2.Range("E10").Value2 = "3-8"
$ws2.Range("F10").Value2 = 15
$ws2.Range("G10").Value2 = "Bangalore"
$ws2.Range("H10").Value2 = "1 Month"
$ws2.Range("I10").Value2 = 40
$ws2.Range("J10").Value2 = "6-14"
$ws2.Range("K10").Value2 = "swapnil.tamhane@cai.io"
$ws2.Rows.Item(10).RowHeight = 29

$ws2.Hyperlinks.Add($ws2.Range("K10"), "mailto:swapnil.tamhane@cai.io", [Type]::Missing, [Type]::Missing, "swapnil.tamhane@cai.io") | Out-Null

# ---------------------------------------------------------------------------
# Step 2: write the new ServiceNow PM posting into NaukriSearch!A2:K2, reusing
# the formatting that is already used by Sheet1!A9:K9 (border style family
# "11/12/13/14" already present in styles.xml)
# ---------------------------------------------------------------------------
$ws2.Range("B9").Copy() | Out-Null
$ws1.Range("A2").PasteSpecial(-4122) | Out-Null
$ws1.Range("F2").PasteSpecial(-4122) | Out-Null
$ws1.Range("G2").PasteSpecial(-4122) | Out-Null
$ws1.Range("H2").PasteSpecial(-4122) | Out-Null
$ws1.Range("I2").PasteSpecial(-4122) | Out-Null
$ws1.Range("B2").PasteSpecial(-4122) | Out-Null

$ws2.Range("C9").Copy() | Out-Null
$ws1.Range("C2").PasteSpecial(-4122) | Out-Null
$ws1.Range("D2").PasteSpecial(-4122) | Out-Null

$ws2.Range("E9").Copy() | Out-Null
$ws1.Range("E2").PasteSpecial(-4122) | Out-Null
$ws1.Range("J2").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Value2 = 1
$ws1.Range("B2").Value2 = "ServiceNow PM"
$ws1.Range("C2").Value2 = "ServiceNow, ITSM, Service Portals"
$ws1.Range("D2").Value2 = "ServiceNow, ITSM, Project Management, Business Analyst, Incident Management, Problem Management, Configuration Management, Change Management, Knowledge Management"
$ws1.Range("E2").Value2 = "6-12"
$ws1.Range("F2").Value2 = 15
$ws1.Range("G2").Value2 = "Bangalore"
$ws1.Range("H2").Value2 = "Immediate"
$ws1.Range("I2").Value2 = 25
$ws1.Range("J2").Value2 = "7-22"
# K2 (email/hyperlink) is unchanged - still swapnil.tamhane@cai.io

$ws1.Rows.Item(2).RowHeight = 42.5
$ws1.Columns.Item(3).ColumnWidth = 17.6

# ---------------------------------------------------------------------------
# Step 3: append the Sitecore posting as Sheet1!A11:J11
# ---------------------------------------------------------------------------
$ws2.Range("A9").Copy() | Out-Null
$ws2.Range("A11").PasteSpecial(-4122) | Out-Null
$ws2.Range("B11").PasteSpecial(-4122) | Out-Null
$ws2.Range("C11").PasteSpecial(-4122) | Out-Null
$ws2.Range("D11").PasteSpecial(-4122) | Out-Null
$ws2.Range("F11").PasteSpecial(-4122) | Out-Null
$ws2.Range("G11").PasteSpecial(-4122) | Out-Null
$ws2.Range("H11").PasteSpecial(-4122) | Out-Null
$ws2.Range("I11").PasteSpecial(-4122) | Out-Null

$ws2.Range("E9").Copy() | Out-Null
$ws2.Range("E11").PasteSpecial(-4122) | Out-Null
$ws2.Range("J11").PasteSpecial(-4122) | Out-Null

$ws2.Range("A11").Value2 = 1
$ws2.Range("B11").Value2 = "Sitecore"
$ws2.Range("C11").Value2 = "Sitecore, Sitecore Certified"
$ws2.Range("D11").Value2 = ".NET, MVC, C#,javascript"
$ws2.Range("E11").Value2 = "3-8"
$ws2.Range("F11").Value2 = "15 Days"
$ws2.Range("G11").Value2 = "Bangalore"
$ws2.Range("H11").Value2 = "1 Month"
$ws2.Range("I11").Value2 = 50
$ws2.Range("J11").Value2 = "12-15"

# ---------------------------------------------------------------------------
# Final selection state (matches the authored workbook)
# ---------------------------------------------------------------------------
$ws2.Range("D9").Select() | Out-Null
$ws1.Range("C2").Select() | Out-Null
